$p = $ppt.ActivePresentation

function Set-ConsolidatedText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    # Force a transient change first so the host always rewrites the
    # paragraph as a single run (no-op text assignments are skipped).
    $tr.Text = "~~~"
    $tr.Text = $text
}

$titles = @{
    1  = "Slide 1 (Content)"
    2  = "Slide 2 (Content)"
    3  = "Slide 3 (Content)"
    4  = "Slide 4 (Content)"
    5  = "Slide 5 (Two Content)"
    6  = "Slide 6 (Two Content Right)"
    7  = "Slide 7 (Content with Caption)"
    8  = "Slide 8 (Comparison)"
    9  = "Slide 9 (Content)"
    10 = "Slide 10 (Content)"
    11 = "Slide 11 (Content)"
    12 = "Slide 12 (Content)"
}

foreach ($slideIdx in $titles.Keys) {
    $s = $p.Slides.Item($slideIdx)
    $titleShape = $s.Shapes.Item(1)
    Set-ConsolidatedText $titleShape $titles[$slideIdx]
}

# Image-caption textboxes that also had their runs split per word.
Set-ConsolidatedText $p.Slides.Item(6).Shapes.Item(3) "an image"
Set-ConsolidatedText $p.Slides.Item(7).Shapes.Item(4) "An image"
Set-ConsolidatedText $p.Slides.Item(8).Shapes.Item(4) "An image"
